$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-3100: add formatter convert to number for some columns
# Columns B..G in data rows (processing, completed, canceled, deferred, closed, new_or_reopened)
# get a ":formatN()" suffix appended to their placeholder text, and switch to a
# numeric ("0") number format. The address column (A) is left untouched.

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in 2..3) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $text = $cell.Value2
        # Insert ":formatN()" right before the closing "}" of the placeholder,
        # e.g. "{d.tickets[i].processing}" -> "{d.tickets[i].processing:formatN()}"
        $newText = $text.Substring(0, $text.Length - 1) + ":formatN()}"
        $cell.Value = $newText
        $cell.NumberFormat = "0"
    }
}
